$wb = $excel.ActiveWorkbook

# Worksheets involved
$wsData = $wb.Worksheets.Item("mercadopago")
$wsWarn = $wb.Worksheets.Item("Evaluation Warning")

# Update the K column sequence numbers (rows 2-6) on the "mercadopago" sheet
$wsData.Range("K2").Value = 47662
$wsData.Range("K3").Value = 47673
$wsData.Range("K4").Value = 47674
$wsData.Range("K5").Value = 47684
$wsData.Range("K6").Value = 47687

# Move the selection on "mercadopago" from K2:K6 to J1 (style tweak for the
# 'Procesar Archivo' button selection state), then restore the original
# active sheet ("Evaluation Warning") so the workbook-level active tab is
# unchanged.
$wsData.Activate()
$wsData.Range("J1").Select()
$wsWarn.Activate()
